$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: Name changes from Burhan -> Akash
$ws.Range("A2").Value = "Akash"

# Row 3: Name changes from Nidhi -> Priya, and Gender (previously shared
# string "Female") becomes an explicit "Female" value again
$ws.Range("A3").Value = "Priya"
$ws.Range("B3").Value = "Female"

# Selection moves from B3 to A3
$ws.Range("A3").Select()
